# Auto-generated Excel COM-interop script
# Applies 2024-04-23 crime-count updates (per commit message) to the
# violent-crime-full-year workbook: 162 cell updates across 43 worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 2219
$ws.Range('K3').Value = 2148
$ws.Range('J4').Value = 1810
$ws.Range('K4').Value = 451
$ws.Range('K6').Value = 2712
$ws.Range('J7').Value = 29279
$ws.Range('K7').Value = 7674

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 153
$ws.Range('K3').Value = 146
$ws.Range('K6').Value = 177
$ws.Range('K7').Value = 515

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 53
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 168

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 111
$ws.Range('K6').Value = 79
$ws.Range('K7').Value = 298

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 43
$ws.Range('K3').Value = 36
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 122

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 87
$ws.Range('K4').Value = 10
$ws.Range('K7').Value = 249

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K6').Value = 42
$ws.Range('K7').Value = 141

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 56
$ws.Range('K5').Value = 13
$ws.Range('K7').Value = 222
$ws.Range('K8').Value = 515
$ws.Range('K11').Value = 165
$ws.Range('K13').Value = 11
$ws.Range('K14').Value = 42
$ws.Range('K15').Value = 77
$ws.Range('K19').Value = 217
$ws.Range('K20').Value = 165
$ws.Range('K21').Value = 22
$ws.Range('K22').Value = 24
$ws.Range('K27').Value = 84
$ws.Range('K29').Value = 385
$ws.Range('K31').Value = 86
$ws.Range('K33').Value = 298
$ws.Range('K36').Value = 93
$ws.Range('K37').Value = 249
$ws.Range('K42').Value = 265
$ws.Range('K43').Value = 71
$ws.Range('K44').Value = 73
$ws.Range('K48').Value = 97
$ws.Range('K50').Value = 50
$ws.Range('K52').Value = 205
$ws.Range('K54').Value = 138
$ws.Range('K55').Value = 83
$ws.Range('K57').Value = 21
$ws.Range('K60').Value = 52
$ws.Range('J63').Value = 97
$ws.Range('K63').Value = 31
$ws.Range('K64').Value = 48
$ws.Range('K67').Value = 295
$ws.Range('K72').Value = 36
$ws.Range('K77').Value = 53
$ws.Range('K78').Value = 104
$ws.Range('K79').Value = 202
$ws.Range('K83').Value = 168
$ws.Range('K85').Value = 375
$ws.Range('K86').Value = 51
$ws.Range('K94').Value = 90
$ws.Range('K95').Value = 122
$ws.Range('K96').Value = 106
$ws.Range('K99').Value = 141
$ws.Range('J101').Value = 29279
$ws.Range('K101').Value = 7674

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 31
$ws.Range('K7').Value = 86

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 88
$ws.Range('K4').Value = 19
$ws.Range('K6').Value = 93
$ws.Range('K7').Value = 295

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K4').Value = 10
$ws.Range('K6').Value = 57
$ws.Range('K7').Value = 138

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 100
$ws.Range('K6').Value = 127
$ws.Range('K7').Value = 385

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K4').Value = 15
$ws.Range('K7').Value = 97

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 69
$ws.Range('K7').Value = 217

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K2').Value = 13
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 73

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K3').Value = 9
$ws.Range('K7').Value = 42

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 64
$ws.Range('K3').Value = 80
$ws.Range('K4').Value = 10
$ws.Range('K6').Value = 109
$ws.Range('K7').Value = 265

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('K4').Value = 2
$ws.Range('K6').Value = 11

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 32
$ws.Range('K3').Value = 25
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 104

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K6').Value = 29
$ws.Range('K7').Value = 83

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 34
$ws.Range('K6').Value = 51
$ws.Range('K7').Value = 106

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('K3').Value = 7
$ws.Range('K7').Value = 22

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 73
$ws.Range('K6').Value = 46
$ws.Range('K7').Value = 202

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K2').Value = 10
$ws.Range('K7').Value = 48

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 55
$ws.Range('K3').Value = 47
$ws.Range('K6').Value = 57
$ws.Range('K7').Value = 165

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 37
$ws.Range('K3').Value = 32
$ws.Range('K7').Value = 93

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K3').Value = 67
$ws.Range('K6').Value = 60
$ws.Range('K7').Value = 222

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K2').Value = 25
$ws.Range('K7').Value = 90

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K6').Value = 26
$ws.Range('K7').Value = 77

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K3').Value = 7
$ws.Range('K6').Value = 30
$ws.Range('K7').Value = 50

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K3').Value = 41
$ws.Range('K6').Value = 67
$ws.Range('K7').Value = 165

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 17
$ws.Range('K7').Value = 56

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('K2').Value = 4
$ws.Range('K7').Value = 13

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 84

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K4').Value = 22
$ws.Range('K7').Value = 51

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K6').Value = 11
$ws.Range('K7').Value = 21

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K6').Value = 15
$ws.Range('K7').Value = 52

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K2').Value = 12
$ws.Range('K7').Value = 71

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 134
$ws.Range('K3').Value = 126
$ws.Range('K4').Value = 20
$ws.Range('K6').Value = 90
$ws.Range('K7').Value = 375

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('K2').Value = 12
$ws.Range('K7').Value = 24

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K6').Value = 20
$ws.Range('K7').Value = 36

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K3').Value = 19
$ws.Range('K7').Value = 53

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 57
$ws.Range('K7').Value = 205
